$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "G2"  = 2.5
    "I2"  = 3.1
    "L2"  = 4
    "Z2"  = 23
    "AA2" = 23
    "BB2" = 351

    "M4"  = 1.13
    "N4"  = 6
    "O4"  = 1.53
    "P4"  = 2.38
    "Q4"  = 2.7
    "R4"  = 1.44
    "S4"  = 1.62
    "T4"  = 2.2
    "W4"  = 8
    "AA4" = 41
    "AI4" = 9
    "AL4" = 23
    "AP4" = 41
    "AT4" = 2.2

    "M5"  = 1.07
    "N5"  = 9

    "G6"  = 1.48
    "H6"  = 4
    "J6"  = 2.05
    "L6"  = 7
    "M6"  = 1.07
    "N6"  = 9
    "O6"  = 1.3
    "P6"  = 3.4
    "R6"  = 1.8
    "S6"  = 1.4
    "T6"  = 2.75
    "U6"  = 2.2
    "V6"  = 1.62
    "Y6"  = 9
    "Z6"  = 9.5
    "AN6" = 3.25
    "AS6" = 201
    "AT6" = 2.75
    "AW6" = 8
    "AX6" = 41
    "BA6" = 201

    "M9"  = 1.06
    "N9"  = 10
    "O9"  = 1.3
    "P9"  = 3.4
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
